$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values for rows 2-5 (column A: transaction ids, column B: cuotas flag)
$ws.Range("A2").Value = 1872782709
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 302618540
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 1002639483
$ws.Range("B4").Value = 0

$ws.Range("A5").Value = 816818518
$ws.Range("B5").Value = 0

# Remove the old extra rows 6-12 so the sheet shrinks back to A1:B5
$ws.Range("A6:B12").ClearContents()
